$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full A2:D16 replacement block in the new row order, then
# write it in one shot so the reordered resource rows land correctly.
$data = @(
    @('pm_1', 204, 3, '/pm/punch_gill, /pm/punch_recesses, /pm/punch_ribbing'),
    @('sm_1', 378, 2, '/sm/sort, /sm/transport'),
    @('vgr_2', 885, 1, '/vgr/pick_up_and_transport'),
    @('ov_1', 612, 2, '/ov/burn, /ov/temper'),
    @('ov_2', 330, 1, '/ov/burn'),
    @('mm_1', 576, 4, '/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to'),
    @('hw_1', 522, 1, '/hw/human_review'),
    @('mm_2', 381, 4, '/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to'),
    @('hbw_1', 873, 4, '/hbw/get_empty_bucket, /hbw/store, /hbw/store_empty_bucket, /hbw/unload'),
    @('sm_2', 309, 2, '/sm/sort, /sm/transport'),
    @('wt_2', 330, 1, '/wt/pick_up_and_transport'),
    @('dm_2', 177, 3, '/dm/cylindrical_drill, /dm/drill, /dm/lower'),
    @('vgr_1', 1866, 1, '/vgr/pick_up_and_transport'),
    @('hbw_2', 1581, 2, '/hbw/store_empty_bucket, /hbw/unload'),
    @('wt_1', 447, 1, '/wt/pick_up_and_transport')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Drop the old "First Subprocess ID" / "Parameter Keys" columns entirely
# (not just their contents) so the used range shrinks to A:D, matching
# the new sensordata-frequency-checking layout.
$ws.Range("E1:F16").EntireColumn.Delete()
